$p = $ppt.ActivePresentation
$s = $p.Slides.Item(7)
$seq = $s.TimeLine.MainSequence
$eff = $seq.AddEffect($s.Shapes.Item(3), 10)
